# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q2" sheet (same 8-column fund-holdings layout/styles)
#    and place it right after "2021-Q2", renamed to "2022-Q1".
# 2. Fill in the 2022-Q1 fund-holdings data (5 rows).
# 3. Update the "总计" (totals) sheet with a new leading row for 2022-Q1,
#    shifting the previous rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by copying "2021-Q2" (keeps header/style)
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")
$template.Copy($null, $template)

$newSheet = $wb.Worksheets.Item("2021-Q2 (2)")
$newSheet.Name = "2022-Q1"

# Extend formatting of column A (index column, bold+border style) down to
# rows 4-6 (template sheet only had rows 2-3 populated).
$newSheet.Range("A3").Copy()
$newSheet.Range("A4:A6").PasteSpecial(-4122)

# The template's D1 header ("基金金额") differs from the new sheet's header
# ("基金规模") - fix it up (keeps the s="2" style already on the cell).
$newSheet.Cells.Item(1,4).Value = "基金规模"

# Helper-less inline pattern for writing "numeric-looking" values as TEXT
# (NumberFormat "@" forces text entry, then Style reset back to Normal
# keeps the default/unstyled look of the data cells while Value stays text).

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 : 001411 / 诺安创新驱动灵活配置混合A
$newSheet.Cells.Item(2,1).Value = 0
Set-TextValue $newSheet.Cells.Item(2,2) "001411"
$newSheet.Cells.Item(2,3).Value = "诺安创新驱动灵活配置混合A"
Set-TextValue $newSheet.Cells.Item(2,4) "3.96"
Set-TextValue $newSheet.Cells.Item(2,5) "94.02"
Set-TextValue $newSheet.Cells.Item(2,6) "4.47"
Set-TextValue $newSheet.Cells.Item(2,7) "0.1770"
$newSheet.Cells.Item(2,8).Value = 7

# Row 3 : 005495 / 创金合信科技成长主题股票A
$newSheet.Cells.Item(3,1).Value = 1
Set-TextValue $newSheet.Cells.Item(3,2) "005495"
$newSheet.Cells.Item(3,3).Value = "创金合信科技成长主题股票A"
Set-TextValue $newSheet.Cells.Item(3,4) "2.62"
Set-TextValue $newSheet.Cells.Item(3,5) "84.91"
Set-TextValue $newSheet.Cells.Item(3,6) "2.75"
Set-TextValue $newSheet.Cells.Item(3,7) "0.0720"
$newSheet.Cells.Item(3,8).Value = 6

# Row 4 : 002051 / 诺安创新驱动灵活配置混合C
$newSheet.Cells.Item(4,1).Value = 2
Set-TextValue $newSheet.Cells.Item(4,2) "002051"
$newSheet.Cells.Item(4,3).Value = "诺安创新驱动灵活配置混合C"
Set-TextValue $newSheet.Cells.Item(4,4) "1.33"
Set-TextValue $newSheet.Cells.Item(4,5) "94.02"
Set-TextValue $newSheet.Cells.Item(4,6) "4.47"
Set-TextValue $newSheet.Cells.Item(4,7) "0.0595"
$newSheet.Cells.Item(4,8).Value = 7

# Row 5 : 005496 / 创金合信科技成长主题股票C
$newSheet.Cells.Item(5,1).Value = 3
Set-TextValue $newSheet.Cells.Item(5,2) "005496"
$newSheet.Cells.Item(5,3).Value = "创金合信科技成长主题股票C"
Set-TextValue $newSheet.Cells.Item(5,4) "0.72"
Set-TextValue $newSheet.Cells.Item(5,5) "84.91"
Set-TextValue $newSheet.Cells.Item(5,6) "2.75"
Set-TextValue $newSheet.Cells.Item(5,7) "0.0198"
$newSheet.Cells.Item(5,8).Value = 6

# Row 6 : 002020 / 国都创新驱动灵活配置混合
$newSheet.Cells.Item(6,1).Value = 4
Set-TextValue $newSheet.Cells.Item(6,2) "002020"
$newSheet.Cells.Item(6,3).Value = "国都创新驱动灵活配置混合"
Set-TextValue $newSheet.Cells.Item(6,4) "0.15"
Set-TextValue $newSheet.Cells.Item(6,5) "74.87"
Set-TextValue $newSheet.Cells.Item(6,6) "2.33"
Set-TextValue $newSheet.Cells.Item(6,7) "0.0035"
$newSheet.Cells.Item(6,8).Value = 7

# ---------------------------------------------------------------------------
# Step 2: update "总计" sheet - insert the new 2022-Q1 row at the top of the
# data (row 2) and shift the previous three rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the styled index column down to the new row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

# Shift existing data rows 2021-Q2 / 2021-Q1 / 2020-Q4 down by one row
# (old row4 -> new row5, old row3 -> new row4, old row2 -> new row3).
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2020-Q4"
$total.Cells.Item(5,3).Value = 4
$total.Cells.Item(5,4).Value = 0.91

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q1"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.11

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q2"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.11

# New row for 2022-Q1.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 5
$total.Cells.Item(2,4).Value = 0.33

# ---------------------------------------------------------------------------
# Restore the originally active sheet/selection (copying a sheet activates
# the new copy as a side effect).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
